$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Settings
$ws2 = $wb.Worksheets.Item(2)   # Constants
$ws3 = $wb.Worksheets.Item(3)   # Assets

# ---------------------------------------------------------------
# Sheet "Settings": rework rows 5-8 (new UiPath project settings)
# ---------------------------------------------------------------

# Row 5: libro1 -> fileSalida (points at the "salida" folder)
$ws1.Range("A5").Value = "fileSalida"
$ws1.Range("B5").Value = "C:\RPA\Uipath_RPA_CristianAzana\salida\"
$ws1.Rows.Item(5).RowHeight = 15

# Row 6: file_salida -> web_URL
$ws1.Range("A6").Value = "web_URL"
$ws1.Range("B6").Value = "https://micnt.com.ec/cntapp/guia104/php/guia_cntat.php?hflagsubmit=0&cmbcriterio=3&cmbprov2=17&txtusuarioapellido=&txtusuarionombre=&captchaSelection="
$ws1.Range("B6").WrapText = $true
$ws1.Range("C6").Value = "url del navegador"

# New row 7: fileExcel
$ws1.Range("A7").Value = "fileExcel"
$ws1.Range("B7").Value = "C:\RPA\Uipath_RPA_CristianAzana\salida\Libro.xlsx"
$ws1.Rows.Item(7).RowHeight = 14.25

# New row 8: fileArchivo
$ws1.Range("A8").Value = "fileArchivo"
$ws1.Range("B8").Value = "C:\RPA\Uipath_RPA_CristianAzana\archivo\"
$ws1.Rows.Item(8).RowHeight = 14.25

# Remove the now-superfluous trailing blank row (949 -> 948 rows used)
$ws1.Rows.Item(949).Delete()

# Column sizing: widen A and B, drop the old "best fit" widths
$ws1.Columns.Item(1).ColumnWidth = 41.166666666666664
$ws1.Columns.Item(2).ColumnWidth = 59.592447916666664

# ---------------------------------------------------------------
# Sheet "Constants": add rerunMaxCount / count rows
# ---------------------------------------------------------------
$ws2.Range("A2").Value = "rerunMaxCount"
$ws2.Range("B2").Value = 3
$ws2.Range("C2").Value = "numero maximo de intentos"

$ws2.Range("A3").Value = "count"
$ws2.Range("B3").Value = 0
$ws2.Range("C3").Value = "inicio del contador"

# ---------------------------------------------------------------
# View state: zoom + selection per sheet (Settings stays active last)
# ---------------------------------------------------------------
$ws3.Activate()
$excel.ActiveWindow.Zoom = 70
$ws3.Range("B24").Select()

$ws2.Activate()
$excel.ActiveWindow.Zoom = 85
$ws2.Range("C19").Select()

$ws1.Activate()
$ws1.Range("B12").Select()
